{"js": "// Add two new bullet items to the \"Tools and Software\" list, right after\n// the existing \"GIMP\" / \"Image creation and editing program for creating\n// 2D art assets\" entry and before the section break:\n//   - \"GitHub\" (tool-name level, ilvl=1)\n//   - \"Version control program to be used when creating the game and\n//     server\" (description level, ilvl=2)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the last paragraph in the document (the GIMP description line);\n// the new bullets land right after it, mirroring its list/run formatting.\nconst items = paragraphs.items;\nconst anchor = items[items.length - 1];\n\nconst rPr =\n  '<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n  '<w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>';\n\nfunction listParaXml(ilvl, text, spellCheck) {\n  const pPr =\n    '<w:pPr><w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"' + ilvl + '\"/><w:numId w:val=\"3\"/></w:numPr>' +\n    '<w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/>' +\n    rPr +\n    '</w:pPr>';\n  const run = '<w:r>' + rPr + '<w:t>' + text + '</w:t></w:r>';\n  const body = spellCheck\n    ? '<w:proofErr w:type=\"spellStart\"/>' + run + '<w:proofErr w:type=\"spellEnd\"/>'\n    : run;\n  return '<w:p>' + pPr + body + '</w:p>';\n}\n\nconst innerXml =\n  listParaXml(1, \"GitHub\", true) +\n  listParaXml(2, \"Version control program to be used when creating the game and server\", false);\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + innerXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// A collapsed range positioned right after the anchor paragraph; inserting\n// OOXML there with InsertLocation.After adds new sibling paragraphs without\n// disturbing the anchor paragraph's own content.\nconst afterRange = anchor.getRange(Word.RangeLocation.after);\nafterRange.insertOoxml(flatOpc, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add two new bullet items to the \"Tools and Software\" list, right after\n# the existing \"GIMP\" / \"Image creation and editing program for creating\n# 2D art assets\" entry and at the very end of the document body:\n#   - \"GitHub\" (tool-name level -> w:ilvl=1, i.e. ListLevelNumber=2)\n#   - \"Version control program to be used when creating the game and\n#     server\" (description level -> w:ilvl=2, i.e. ListLevelNumber=3)\n\n$d = $word.ActiveDocument\n\n# Collapsed range at the very end of the document body content (just before\n# the final section mark) \u2014 inserting WordML there appends new sibling\n# paragraphs without disturbing the existing last paragraph.\n$endRange = $d.Range($d.Content.End, $d.Content.End)\n\n$xml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"3\"/></w:numPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>GitHub</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"2\"/><w:numId w:val=\"3\"/></w:numPr><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Version control program to be used when creating the game and server</w:t></w:r></w:p>\n'@\n\n$endRange.InsertXML($xml)\n"}
